# Refresh cryptocurrency Price (D) and Volume(1h) (E) columns with latest scraped values.
# Some "Price" strings are numeric-looking (e.g. "310.77") but must remain as literal
# text cells (matching the original inlineStr cells), so we temporarily mark the
# range as Text before assigning, then restore the original cell style/format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $originalStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $originalStyle
}

$ws.Range("D2").Value = "27.042.83"
$ws.Range("E2").Value = "  +2.19%  "
$ws.Range("D3").Value = "1.847.20"
$ws.Range("E3").Value = "  +2.29%  "
Set-TextValue $ws.Range("D4") "1.006"
$ws.Range("E4").Value = "  +0.14%  "
Set-TextValue $ws.Range("D5") "310.77"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("E6").Value = "  +0.20%  "
Set-TextValue $ws.Range("D7") "0.4693"
$ws.Range("E7").Value = "  +3.83%  "
Set-TextValue $ws.Range("D8") "0.3640"
$ws.Range("E8").Value = "  +1.18%  "
Set-TextValue $ws.Range("D9") "0.07185"
$ws.Range("E9").Value = "  +1.60%  "
Set-TextValue $ws.Range("D10") "0.9393"
$ws.Range("E10").Value = "  +5.69%  "
Set-TextValue $ws.Range("D11") "19.67"
$ws.Range("E11").Value = "  +1.20%  "
Set-TextValue $ws.Range("D12") "0.07684"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("D13").Value = "1.862.22"
$ws.Range("E13").Value = "  +5.49%  "
Set-TextValue $ws.Range("D14") "5.299"
$ws.Range("E14").Value = "  +0.16%  "
Set-TextValue $ws.Range("D15") "6.397"
$ws.Range("E15").Value = "  +1.17%  "
Set-TextValue $ws.Range("D16") "88.56"
$ws.Range("E16").Value = "  +3.89%  "
Set-TextValue $ws.Range("D17") "1.007"
$ws.Range("E17").Value = "  +0.06%  "
Set-TextValue $ws.Range("D18") "0.000008600"
$ws.Range("E18").Value = "  +1.37%  "
Set-TextValue $ws.Range("D19") "1.008"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").Value = "27.039.82"
$ws.Range("E20").Value = "  +2.04%  "
Set-TextValue $ws.Range("D21") "14.37"
$ws.Range("E21").Value = "  +1.25%  "
Set-TextValue $ws.Range("D22") "5.042"
$ws.Range("E22").Value = "  +1.56%  "
Set-TextValue $ws.Range("D23") "10.65"
$ws.Range("E23").Value = "  +1.29%  "
$ws.Range("E24").Value = "  -1.81%  "
Set-TextValue $ws.Range("D25") "152.41"
$ws.Range("E25").Value = "  +0.74%  "
Set-TextValue $ws.Range("D26") "18.04"
$ws.Range("E26").Value = "  +1.33%  "
Set-TextValue $ws.Range("D27") "2.033"
$ws.Range("E27").Value = "  -1.66%  "
Set-TextValue $ws.Range("D28") "114.18"
$ws.Range("E28").Value = "  +1.87%  "
Set-TextValue $ws.Range("D29") "4.936"
$ws.Range("E29").Value = "  +1.56%  "
Set-TextValue $ws.Range("D30") "0.08853"
$ws.Range("E30").Value = "  +1.78%  "
$ws.Range("E31").Value = "  +2.18%  "
Set-TextValue $ws.Range("D32") "2.852"
$ws.Range("E32").Value = "  +0.76%  "
Set-TextValue $ws.Range("D33") "1.187"
$ws.Range("E33").Value = "  +7.66%  "
Set-TextValue $ws.Range("D34") "0.7502"
$ws.Range("E34").Value = "  +4.06%  "
Set-TextValue $ws.Range("D35") "4.478"
$ws.Range("E35").Value = "  +0.55%  "
Set-TextValue $ws.Range("D36") "1.088"
$ws.Range("E36").Value = "  +1.29%  "
Set-TextValue $ws.Range("D37") "2.982"
$ws.Range("E37").Value = "  +3.14%  "
Set-TextValue $ws.Range("D38") "0.01941"
$ws.Range("E38").Value = "  +0.55%  "
Set-TextValue $ws.Range("D39") "0.05164"
$ws.Range("E39").Value = "  +1.29%  "
Set-TextValue $ws.Range("D40") "0.5152"
$ws.Range("E40").Value = "  +0.72%  "
Set-TextValue $ws.Range("D41") "6.937"
$ws.Range("E41").Value = "  +2.30%  "
Set-TextValue $ws.Range("D42") "0.1516"
$ws.Range("E42").Value = "  +0.37%  "
Set-TextValue $ws.Range("D43") "8.202"
$ws.Range("E43").Value = "  +2.35%  "
Set-TextValue $ws.Range("D44") "0.4729"
$ws.Range("E44").Value = "  +1.45%  "
Set-TextValue $ws.Range("D45") "10.34"
$ws.Range("E45").Value = "  +3.13%  "
Set-TextValue $ws.Range("D46") "1.007"
$ws.Range("E46").Value = "  +0.32%  "
Set-TextValue $ws.Range("D47") "100.30"
$ws.Range("E47").Value = "  -0.17%  "
Set-TextValue $ws.Range("D48") "1.605"
$ws.Range("E48").Value = "  +2.02%  "
Set-TextValue $ws.Range("D49") "0.06058"
$ws.Range("E49").Value = "  +1.34%  "
Set-TextValue $ws.Range("D50") "64.23"
$ws.Range("E50").Value = "  +0.62%  "
Set-TextValue $ws.Range("D51") "36.23"
$ws.Range("E51").Value = "  +0.24%  "
